# Update values in Sheet1 (result_data_RandomForest) per "Update Name of Algo" commit.
# These values correspond to the RandomForest imputation results being regenerated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 16.2365
$ws.Range("C7").Value = -13.0462
$ws.Range("B8").Value = 5.886199999999999
$ws.Range("B10").Value = 5.3579
$ws.Range("B12").Value = 4.7867
$ws.Range("C15").Value = -13.74069999999999
$ws.Range("B18").Value = 7.469199999999995
$ws.Range("C18").Value = -11.57409999999999
$ws.Range("E18").Value = 18.07700000000001
$ws.Range("E19").Value = 16.43850000000001
$ws.Range("C20").Value = -11.9377
$ws.Range("E27").Value = 16.56519999999999
$ws.Range("C29").Value = -11.43570000000001
$ws.Range("C30").Value = -12.91779999999999
$ws.Range("C31").Value = -13.14349999999999
$ws.Range("E31").Value = 16.5398
$ws.Range("B37").Value = 8.830000000000002
$ws.Range("E38").Value = 16.385
$ws.Range("C40").Value = -12.90530000000001
$ws.Range("E42").Value = 16.39229999999999
$ws.Range("E44").Value = 16.6071
$ws.Range("E47").Value = 16.68759999999999
$ws.Range("C50").Value = -13.05629999999999
$ws.Range("B55").Value = 6.601199999999997
$ws.Range("E58").Value = 16.14270000000002
$ws.Range("E65").Value = 17.34260000000001
$ws.Range("B68").Value = 6.1562
$ws.Range("C68").Value = -12.08690000000001
$ws.Range("E73").Value = 17.41210000000002
$ws.Range("C76").Value = -12.4139
$ws.Range("B77").Value = 9.152800000000006
$ws.Range("B78").Value = 9.371200000000004
$ws.Range("B81").Value = 5.276600000000005
$ws.Range("B82").Value = 5.212100000000001
$ws.Range("C87").Value = -13.90469999999999
$ws.Range("C88").Value = -12.79539999999999
$ws.Range("E90").Value = 16.46069999999999
$ws.Range("E94").Value = 18.98490000000002
$ws.Range("E95").Value = 18.17660000000002
$ws.Range("C96").Value = -12.79950000000001
$ws.Range("C98").Value = -12.12229999999999
$ws.Range("C101").Value = -12.16080000000001
$ws.Range("E101").Value = 16.6445
$ws.Range("C102").Value = -13.06810000000001
